$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("stats")

$ws.Range("D2").Value = 0.0001352312974631786
$ws.Range("E2").Value = 0.03035148791968822
$ws.Range("G2").Value = 0.002112641464918852
$ws.Range("H2").Value = 0.004344364628195763
$ws.Range("I2").Value = 0.008957847021520138
$ws.Range("J2").Value = 0.01196514163166285
$ws.Range("K2").Value = 0.0008195945993065834
$ws.Range("D3").Value = 0.00185466418042779
$ws.Range("E3").Value = 0.03206227114424109
$ws.Range("G3").Value = 0.002093623857945204
$ws.Range("H3").Value = 0.005747798830270767
$ws.Range("I3").Value = 0.008508458267897367
$ws.Range("J3").Value = 0.01302590081468225
$ws.Range("K3").Value = 0.0007318048737943172
$ws.Range("D4").Value = 0.001227468717843294
$ws.Range("E4").Value = 0.02108264388516545
$ws.Range("G4").Value = 0.001259608194231987
$ws.Range("H4").Value = 0.003818481229245663
$ws.Range("I4").Value = 0.005760840140283108
$ws.Range("J4").Value = 0.008561976719647646
$ws.Range("K4").Value = 0.0005278461612761021
$ws.Range("D5").Value = 0.0001407810486853123
$ws.Range("E5").Value = 0.01870762603357434
$ws.Range("G5").Value = 0.001327746547758579
$ws.Range("H5").Value = 0.0025911470875144
$ws.Range("I5").Value = 0.005706444848328829
$ws.Range("J5").Value = 0.007380329538136721
$ws.Range("K5").Value = 0.0005363314412534237
$ws.Range("D6").Value = 0.002037713769823313
$ws.Range("E6").Value = 0.07072611898183823
$ws.Range("G6").Value = 0.002854987047612667
$ws.Range("H6").Value = 0.008526691701263189
$ws.Range("I6").Value = 0.0453035244718194
$ws.Range("J6").Value = 0.01018311083316803
$ws.Range("K6").Value = 0.001093864440917969
$ws.Range("D8").Value = 0.0001352312974631786
$ws.Range("E8").Value = 0.03035148791968822
$ws.Range("G8").Value = 0.002112641464918852
$ws.Range("H8").Value = 0.004344364628195763
$ws.Range("I8").Value = 0.008957847021520138
$ws.Range("J8").Value = 0.01196514163166285
$ws.Range("K8").Value = 0.0008195945993065834
$ws.Range("D9").Value = 0.00185466418042779
$ws.Range("E9").Value = 0.03206227114424109
$ws.Range("G9").Value = 0.002093623857945204
$ws.Range("H9").Value = 0.005747798830270767
$ws.Range("I9").Value = 0.008508458267897367
$ws.Range("J9").Value = 0.01302590081468225
$ws.Range("K9").Value = 0.0007318048737943172
$ws.Range("D10").Value = 0.001227468717843294
$ws.Range("E10").Value = 0.02108264388516545
$ws.Range("G10").Value = 0.001259608194231987
$ws.Range("H10").Value = 0.003818481229245663
$ws.Range("I10").Value = 0.005760840140283108
$ws.Range("J10").Value = 0.008561976719647646
$ws.Range("K10").Value = 0.0005278461612761021
$ws.Range("D11").Value = 0.0001407810486853123
$ws.Range("E11").Value = 0.01870762603357434
$ws.Range("G11").Value = 0.001327746547758579
$ws.Range("H11").Value = 0.0025911470875144
$ws.Range("I11").Value = 0.005706444848328829
$ws.Range("J11").Value = 0.007380329538136721
$ws.Range("K11").Value = 0.0005363314412534237
$ws.Range("D12").Value = 0.002037713769823313
$ws.Range("E12").Value = 0.07072611898183823
$ws.Range("G12").Value = 0.002854987047612667
$ws.Range("H12").Value = 0.008526691701263189
$ws.Range("I12").Value = 0.0453035244718194
$ws.Range("J12").Value = 0.01018311083316803
$ws.Range("K12").Value = 0.001093864440917969
